$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CSS-Accept-DB")
$ws2 = $wb.Worksheets.Item("CSS-Reject-DB")

# --- Sheet "CSS-Accept-DB" (sheet1) ---
# K6: expirationDate 2021-09-21 -> 2021-11-21
$ws1.Range("K6").Value = "{`n    ""status"": ""QUOTE_RECEIVED"",`n    ""expirationDate"": ""2021-11-21T04:59:00.000Z"",`n    ""insurancePremium"": {`n        ""amount"": 500,`n        ""currency"": ""CHF""`n    },`n    ""policyLimit"": {`n        ""amount"": 50000,`n        ""currency"": ""CHF""`n    }`n}"

# Row 8: DB verify step - G8 VERIFY -> SELECT, K8 gets the select query text
$ws1.Range("G8").Value = "SELECT"
$ws1.Range("K8").Value = "select iqr.id, iq.insurance_premium_amount, iq.insurance_premium_currency, iq.policy_limit_amount from insurancequotes iq INNER JOIN insurancequoterequests iqr on iq.id = iqr.insurance_quote_id and iqr.id  =  [quoteId]"

# --- Sheet "CSS-Reject-DB" (sheet2) ---
# J6: expirationDate 2021-09-21 -> 2022-09-21
$ws2.Range("J6").Value = "{`n    ""status"": ""QUOTE_RECEIVED"",`n    ""expirationDate"": ""2022-09-21T04:59:00.000Z"",`n    ""insurancePremium"": {`n        ""amount"": 500,`n        ""currency"": ""CHF""`n    },`n    ""policyLimit"": {`n        ""amount"": 50000,`n        ""currency"": ""CHF""`n    }`n}"

# Row 8: DB verify step restructured
$ws2.Range("E8").Value = "Read Quote information"
$ws2.Range("F8").Value = "SELECT"
$ws2.Range("G8").Value = ""
$ws2.Range("J8").Value = "select iqr.id, iq.insurance_premium_amount, iq.insurance_premium_currency, iq.policy_limit_amount from insurancequotes iq INNER JOIN insurancequoterequests iqr on iq.id = iqr.insurance_quote_id and iqr.id  =  [rejectQuoteId]"
$ws2.Range("L8").Value = "policy_limit_amount=[0].policy_limit_amount"
$ws2.Range("M8").Value = ""
$ws2.Range("N8").Value = ""
